$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 - copy header style from an existing header cell (B1) then set values
$ws.Range("B1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-64 for columns I (I0) and J (IF)
$data = @(
    @(3, 4),
    @(6, 7),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(5, 6),
    @(6, 6),
    @(6, 7),
    @(6, 7),
    @(7, 7),
    @(7, 7),
    @(4, 5),
    @(7, 7),
    @(8, 8),
    @(6, 7),
    @(7, 7),
    @(6, 6),
    @(9, 9),
    @(7, 8),
    @(6, 7),
    @(5, 6),
    @(7, 7),
    @(5, 6),
    @(6, 6),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(7, 7),
    @(6, 6),
    @(7, 8),
    @(11, 12),
    @(8, 8),
    @(6, 7),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(5, 7),
    @(6, 7),
    @(7, 7),
    @(7, 7),
    @(7, 8),
    @(7, 8),
    @(7, 7),
    @(7, 7),
    @(8, 9),
    @(6, 6),
    @(6, 7),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(6, 7),
    @(7, 7),
    @(5, 6),
    @(6, 6),
    @(6, 7),
    @(6, 7),
    @(6, 6),
    @(6, 8),
    @(6, 6),
    @(6, 7),
    @(5, 6),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
